$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a fifth column (E) of minutes data for the 16/08/2023 meeting ---
# Mirror formatting from column D (the previous week's column) into column E,
# then overwrite the values/text that differ for the new week.

# E3: meeting date (16/08/2023 -> Excel serial 45154), same date format as D3
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = 45154

# E5: Participants present - same text as C5/D5
$ws.Range("D5").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = $ws.Range("D5").Value2

# E6: Participants absent - stays blank, just carry the formatting over
$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial(-4122)

# E8: Subjects discussed - new text for this week
$ws.Range("D8").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value = "Produção na aula de PI"

# E9: Closing notes - stays blank, just carry the formatting over
$ws.Range("D9").Copy()
$ws.Range("E9").PasteSpecial(-4122)

# Give column E its own width (close to the other task columns)
$ws.Columns.Item(5).ColumnWidth = 29.65

# --- New textbox summarising what happened in this week's class ---
$shp = $ws.Shapes.AddTextbox(1, 320.6632283464567, 214.79590551181101, 173.80102362204724, 123.46937007874016)
$shp.Name = "CaixaDeTexto 3"
$shp.TextFrame.Characters().Text = "Nos juntamos em grupo na aula. Separamos tarefas para tofods do grupo, onde alguns codaram e outros fizeram a proto-persona e o mind-Map . Conversamos sobre possiveis itens na area de usuário do protótipo"

# Move the active selection like the saved workbook (cell F9)
$ws.Range("F9").Select()
